# Updated cryptos list on Thu Jan 18 16:43:34 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$ref, [string]$val) {
    # Force the cell to Text so purely-numeric-looking strings (e.g. "34.60")
    # keep their exact literal formatting instead of being auto-coerced into a
    # Number by Excel (which would silently drop trailing zeros, etc).
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    # Drop back to the default style so we do not leave a stray NumberFormat/
    # quotePrefix behind on a cell that had no explicit style originally.
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '42.178.18'
Set-TextCell 'E2' '  -0.41%  '
Set-TextCell 'D3' '2.491.66'
Set-TextCell 'E3' '  -1.77%  '
Set-TextCell 'E4' '  +0.10%  '
Set-TextCell 'D5' '312.55'
Set-TextCell 'E5' '  +0.22%  '
Set-TextCell 'D6' '96.23'
Set-TextCell 'E6' '  -2.77%  '
Set-TextCell 'D7' '0.559'
Set-TextCell 'E7' '  -1.57%  '
Set-TextCell 'E8' '  +0.03%  '
Set-TextCell 'D9' '0.515'
Set-TextCell 'E9' '  -2.17%  '
Set-TextCell 'D10' '34.60'
Set-TextCell 'E10' '  -2.68%  '
Set-TextCell 'D11' '0.0791'
Set-TextCell 'E11' '  -1.42%  '
Set-TextCell 'D12' '0.109'
Set-TextCell 'E12' '  +1.08%  '
Set-TextCell 'D13' '7.10'
Set-TextCell 'E13' '  -3.30%  '
Set-TextCell 'D14' '2.893.94'
Set-TextCell 'E14' '  -1.23%  '
Set-TextCell 'D15' '2.511.40'
Set-TextCell 'E15' '  -3.29%  '
Set-TextCell 'D16' '14.91'
Set-TextCell 'E16' '  -6.43%  '
Set-TextCell 'D17' '0.794'
Set-TextCell 'E17' '  -4.90%  '
Set-TextCell 'D18' '42.251.46'
Set-TextCell 'E18' '  -0.28%  '
Set-TextCell 'D19' '6.45'
Set-TextCell 'E19' '  -4.90%  '
Set-TextCell 'D20' '0.0₃0927'
Set-TextCell 'E20' '  -1.98%  '
Set-TextCell 'D21' '11.88'
Set-TextCell 'E21' '  -2.31%  '
Set-TextCell 'D22' '69.16'
Set-TextCell 'E22' '  +0.59%  '
Set-TextCell 'D23' '238.97'
Set-TextCell 'E23' '  -1.61%  '
Set-TextCell 'D24' '2.82'
Set-TextCell 'E24' '  -2.35%  '
Set-TextCell 'D25' '1.95'
Set-TextCell 'E25' '  -4.32%  '
Set-TextCell 'D26' '1.00'
Set-TextCell 'E26' '  -0.09%  '
Set-TextCell 'D27' '24.99'
Set-TextCell 'E27' '  -5.01%  '
Set-TextCell 'D28' '2.24'
Set-TextCell 'E28' '  -4.46%  '
Set-TextCell 'D29' '9.84'
Set-TextCell 'E29' '  -2.62%  '
Set-TextCell 'D30' '37.21'
Set-TextCell 'E30' '  -6.92%  '
Set-TextCell 'B31' 'Monero'
Set-TextCell 'C31' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D31' '155.43'
Set-TextCell 'E31' '  -1.56%  '
Set-TextCell 'B32' 'Filecoin'
Set-TextCell 'C32' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D32' '5.74'
Set-TextCell 'E32' '  +0.53%  '
Set-TextCell 'D33' '2.66'
Set-TextCell 'E33' '  -6.28%  '
Set-TextCell 'E34' '  +0.45%  '
Set-TextCell 'D35' '0.0771'
Set-TextCell 'E35' '  -3.30%  '
Set-TextCell 'D36' '3.04'
Set-TextCell 'E36' '  -3.46%  '
Set-TextCell 'D37' '1.92'
Set-TextCell 'E37' '  -5.32%  '
Set-TextCell 'E38' '  -3.83%  '
Set-TextCell 'D39' '0.107'
Set-TextCell 'E39' '  -3.47%  '
Set-TextCell 'D40' '0.115'
Set-TextCell 'E40' '  -1.95%  '
Set-TextCell 'D41' '4.10'
Set-TextCell 'E41' '  -2.82%  '
Set-TextCell 'E42' '  -2.85%  '
Set-TextCell 'E43' '  +0.05%  '
Set-TextCell 'D44' '2.010.70'
Set-TextCell 'E44' '  +2.97%  '
Set-TextCell 'D45' '0.0290'
Set-TextCell 'E45' '  -1.98%  '
Set-TextCell 'D46' '3.14'
Set-TextCell 'E46' '  -4.50%  '
Set-TextCell 'D47' '8.73'
Set-TextCell 'E47' '  -1.93%  '
Set-TextCell 'D48' '2.748.66'
Set-TextCell 'E48' '  -1.25%  '
Set-TextCell 'B49' 'BitcoinSV'
Set-TextCell 'C49' 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextCell 'D49' '78.24'
Set-TextCell 'E49' '  -3.07%  '
Set-TextCell 'B50' 'ordi'
Set-TextCell 'C50' 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextCell 'D50' '71.36'
Set-TextCell 'E50' '  -0.92%  '
Set-TextCell 'D51' '0.184'
Set-TextCell 'E51' '  -3.60%  '
